{"js": "// Updated PA3 data values (January 2025) across Table 2.3\n// Each entry is an exact-match, case-sensitive text replacement that is\n// unique within the document body (verified against the source XML).\nconst replacements = [\n  [\"N = 80,139\", \"N = 80,096\"],\n  [\"79,855 (100)\", \"79,812 (100)\"],\n  [\"79,568 (99)\", \"79,525 (99)\"],\n  [\"20,116 (25)\", \"20,101 (25)\"],\n  [\"20,111 (25)\", \"20,104 (25)\"],\n  [\"(464,642]\", \"(464,641]\"],\n  [\"19,999 (25)\", \"19,991 (25)\"],\n  [\"(642,2.39e+03]\", \"(641,2.39e+03]\"],\n  [\"19,913 (25)\", \"19,900 (25)\"],\n  [\"19,998 (25)\", \"19,994 (25)\"],\n  [\"20,207 (25)\", \"20,142 (25)\"],\n  [\"20,010 (25)\", \"20,059 (25)\"],\n  [\"19,924 (25)\", \"19,901 (25)\"],\n  [\"20,219 (25)\", \"20,210 (25)\"],\n  [\"20,067 (25)\", \"20,058 (25)\"],\n  [\"20,003 (25)\", \"19,991 (25)\"],\n  [\"(853,3.37e+03]\", \"(853,2.49e+03]\"],\n  [\"19,850 (25)\", \"19,837 (25)\"],\n  [\"462.8 (319.0, 640.2)\", \"462.8 (318.8, 639.8)\"],\n  [\"707.8 (515.8, 936.6)\", \"707.6 (515.6, 936.2)\"],\n  [\"304.8 (198.3, 443.5)\", \"304.7 (198.3, 443.2)\"],\n  [\"233.3 (115.0, 403.0)\", \"233.0 (115.0, 402.5)\"],\n  [\"583.7 (370.1, 849.9)\", \"583.4 (370.0, 849.4)\"],\n  [\"1,164 (1.5)\", \"1,163 (1.5)\"],\n  [\"45,988 (57)\", \"45,954 (57)\"],\n  [\"34,151 (43)\", \"34,142 (43)\"],\n  [\"6,332 (7.9)\", \"6,319 (7.9)\"],\n  [\"20,164 (25)\", \"20,152 (25)\"],\n  [\"18,921 (24)\", \"18,911 (24)\"],\n  [\"34,722 (43)\", \"34,714 (43)\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found, cannot apply replacement: \"${oldText}\"`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Updated PA3 data values (January 2025) across Table 2.3.\n# Each pair is an exact-match, case-sensitive text replacement that is\n# unique within the document body (verified against the source XML).\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"N = 80,139\", \"N = 80,096\"),\n    @(\"79,855 (100)\", \"79,812 (100)\"),\n    @(\"79,568 (99)\", \"79,525 (99)\"),\n    @(\"20,116 (25)\", \"20,101 (25)\"),\n    @(\"20,111 (25)\", \"20,104 (25)\"),\n    @(\"(464,642]\", \"(464,641]\"),\n    @(\"19,999 (25)\", \"19,991 (25)\"),\n    @(\"(642,2.39e+03]\", \"(641,2.39e+03]\"),\n    @(\"19,913 (25)\", \"19,900 (25)\"),\n    @(\"19,998 (25)\", \"19,994 (25)\"),\n    @(\"20,207 (25)\", \"20,142 (25)\"),\n    @(\"20,010 (25)\", \"20,059 (25)\"),\n    @(\"19,924 (25)\", \"19,901 (25)\"),\n    @(\"20,219 (25)\", \"20,210 (25)\"),\n    @(\"20,067 (25)\", \"20,058 (25)\"),\n    @(\"20,003 (25)\", \"19,991 (25)\"),\n    @(\"(853,3.37e+03]\", \"(853,2.49e+03]\"),\n    @(\"19,850 (25)\", \"19,837 (25)\"),\n    @(\"462.8 (319.0, 640.2)\", \"462.8 (318.8, 639.8)\"),\n    @(\"707.8 (515.8, 936.6)\", \"707.6 (515.6, 936.2)\"),\n    @(\"304.8 (198.3, 443.5)\", \"304.7 (198.3, 443.2)\"),\n    @(\"233.3 (115.0, 403.0)\", \"233.0 (115.0, 402.5)\"),\n    @(\"583.7 (370.1, 849.9)\", \"583.4 (370.0, 849.4)\"),\n    @(\"1,164 (1.5)\", \"1,163 (1.5)\"),\n    @(\"45,988 (57)\", \"45,954 (57)\"),\n    @(\"34,151 (43)\", \"34,142 (43)\"),\n    @(\"6,332 (7.9)\", \"6,319 (7.9)\"),\n    @(\"20,164 (25)\", \"20,152 (25)\"),\n    @(\"18,921 (24)\", \"18,911 (24)\"),\n    @(\"34,722 (43)\", \"34,714 (43)\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n\n    $found = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n    if (-not $found) {\n        Write-Output \"WARNING: text not found, replacement not applied: $oldText\"\n    }\n}\n"}
